$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data driven testing for LoginTest: set both username and password fields to "tutorial"
$ws.Range("A1").Value = "tutorial"
$ws.Range("B1").Value = "tutorial"

# Move the active selection to D6 (from D5)
$ws.Range("D6").Select()
